# Auto-applied edit matching the commit diff for 杭州-漫展信息.xlsx
# Updates view-counts (column F) and refreshes rows 41-43 on the
# '全部类型' sheet to reflect newly scraped event listings.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item('展览')
$wsExhibit.Range("F2").Value = 821
$wsExhibit.Range("F3").Value = 561
$wsExhibit.Range("F4").Value = 308
$wsExhibit.Range("F6").Value = 1157
$wsExhibit.Range("F7").Value = 336
$wsExhibit.Range("F11").Value = 1207
$wsExhibit.Range("F14").Value = 895
$wsExhibit.Range("F15").Value = 881
$wsExhibit.Range("F19").Value = 662
$wsExhibit.Range("F20").Value = 782
$wsExhibit.Range("F21").Value = 1744
$wsExhibit.Range("F22").Value = 2949
$wsExhibit.Range("F23").Value = 857
$wsExhibit.Range("F25").Value = 2227
$wsExhibit.Range("F26").Value = 669
$wsExhibit.Range("F27").Value = 3074
$wsExhibit.Range("F28").Value = 602
$wsExhibit.Range("F29").Value = 359
$wsExhibit.Range("F30").Value = 15
$wsExhibit.Range("F36").Value = 1090
$wsExhibit.Range("F37").Value = 1781
$wsExhibit.Range("F38").Value = 393
$wsExhibit.Range("F41").Value = 196
$wsExhibit.Range("F43").Value = 180
$wsShow = $wb.Worksheets.Item('演出')
$wsShow.Range("F8").Value = 14
$wsAll = $wb.Worksheets.Item('全部类型')
$wsAll.Range("F2").Value = 821
$wsAll.Range("F3").Value = 561
$wsAll.Range("F4").Value = 308
$wsAll.Range("F6").Value = 1157
$wsAll.Range("F7").Value = 336
$wsAll.Range("F10").Value = 1207
$wsAll.Range("F12").Value = 895
$wsAll.Range("F13").Value = 881
$wsAll.Range("F19").Value = 782
$wsAll.Range("F20").Value = 1744
$wsAll.Range("F21").Value = 2949
$wsAll.Range("F22").Value = 857
$wsAll.Range("F25").Value = 2227
$wsAll.Range("F26").Value = 3074
$wsAll.Range("F27").Value = 602
$wsAll.Range("F28").Value = 359
$wsAll.Range("F30").Value = 15
$wsAll.Range("F31").Value = 14
$wsAll.Range("B41").NumberFormat = "@"
$wsAll.Range("B41").Value = '2024-06-08'
$wsAll.Range("B41").ClearFormats()
$wsAll.Range("C41").Value = '杭州·第八届YH樱花动漫游戏文化节'
$wsAll.Range("D41").Value = '德胜东路2539号 梦马汽车小镇'
$wsAll.Range("E41").Value = '2024.06.08 10:00 - 06.10 17:00'
$wsAll.Range("F41").Value = 1090
$wsAll.Range("G41").Value = 65
$wsAll.Range("H41").Value = 'https://show.bilibili.com/platform/detail.html?id=82687'
$wsAll.Range("I41").Value = '//i2.hdslb.com/bfs/openplatform/202403/S5pnadXj1710210939138.png'
$wsAll.Range("B42").NumberFormat = "@"
$wsAll.Range("B42").Value = '2024-06-09'
$wsAll.Range("B42").ClearFormats()
$wsAll.Range("C42").Value = '杭州·第三届日夜国乙only'
$wsAll.Range("D42").Value = '创意路1号 中国智谷富春园区'
$wsAll.Range("E42").Value = '2024.06.09 10:00 - 06.09 23:00'
$wsAll.Range("F42").Value = 1781
$wsAll.Range("G42").Value = 58
$wsAll.Range("H42").Value = 'https://show.bilibili.com/platform/detail.html?id=82618'
$wsAll.Range("I42").Value = '//i2.hdslb.com/bfs/openplatform/202403/fXRzYEFH1710124366279.png'
$wsAll.Range("B43").NumberFormat = "@"
$wsAll.Range("B43").Value = '2024-07-06'
$wsAll.Range("B43").ClearFormats()
$wsAll.Range("C43").Value = '杭州·《爱·永恒》理查德·克莱德曼钢琴音乐会'
$wsAll.Range("D43").Value = '杭州市西湖区省府路9号 浙江省人民大会堂'
$wsAll.Range("E43").Value = '2024.07.06 19:30 - 07.06 21:00'
$wsAll.Range("F43").Value = 1
$wsAll.Range("G43").Value = 980
$wsAll.Range("H43").Value = 'https://show.bilibili.com/platform/detail.html?id=83948'
$wsAll.Range("I43").Value = '//i1.hdslb.com/bfs/openplatform/202404/huCPMql51712639746482.jpeg'
$wsAll.Range("F44").Value = 393
$wsAll.Range("F46").Value = 196
$wsAll.Range("F48").Value = 180
